# [UPDATE] Budget, Report, and Approval
#
# The "rkap" sheet contains budget rows grouped by category (column B,
# driven by shared strings "TA" / "TD" / "TX"). The budget year recorded
# in column A is being rolled forward for the "TA" (rows 2-17) and "TD"
# (rows 18-33) groups, while the "TX" group (rows 34-49) is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rkap")
$ws.Activate()

# TA rows: 2021 -> 2023
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = 2023
}

# TD rows: 2021 -> 2022
for ($r = 18; $r -le 33; $r++) {
    $ws.Cells.Item($r, 1).Value = 2022
}

# Update the view: scroll the sheet so row 35 is at the top and the
# selected cell is C44, matching where the author left off editing.
$win = $excel.ActiveWindow
try { $win.ScrollRow = 35 } catch {}
try { $win.ScrollColumn = 1 } catch {}

$ws.Range("C44").Select()

# Restore the workbook window geometry recorded by Excel on save.
try {
    $win.Left = 12630
    $win.Top = 3000
    $win.Width = 14205
    $win.Height = 11295
} catch {}
